# "started on set MSPeakLists"
# Adds a new "mslists" worksheet (for the upcoming MSPeakLists set-method
# implementation table) after the existing "fGroups" sheet, mirroring its
# legend layout, and moves the active selection over to the new sheet.

$wb = $excel.ActiveWorkbook
$fGroups = $wb.Worksheets.Item(1)

# --- fGroups: selection moves off G35 onto the legend header row ---------
$fGroups.Range("B1:G1").Select() | Out-Null

# --- add the new "mslists" worksheet, placed right after fGroups ---------
$mslists = $wb.Worksheets.Add($null, $fGroups)
$mslists.Name = "mslists"

# Legend header row (as-is / almost as-is / implement / not supported / ionize / done)
$mslists.Range("B1").Value = "as-is"
$mslists.Range("C1").Value = "almost as-is"
$mslists.Range("D1").Value = "implement"
$mslists.Range("E1").Value = "not supported"
$mslists.Range("F1").Value = "ionize"
$mslists.Range("G1").Value = "done"

# Method rows tracking MSPeakLists implementation progress
$mslists.Range("A2").Value = "$"
$mslists.Range("B2").Value = "X"

$mslists.Range("A3").Value = "["
$mslists.Range("C3").Value = "X"

$mslists.Range("A4").Value = "[["
$mslists.Range("C4").Value = "X"

$mslists.Range("A5").Value = "analyses"
$mslists.Range("B5").Value = "X"

$mslists.Range("A6").Value = "as.data.table"
$mslists.Range("C6").Value = "X"

$mslists.Range("A7").Value = "averagedPeakLists"
$mslists.Range("C7").Value = "X"

$mslists.Range("A8").Value = "compoundViewer"
$mslists.Range("E8").Value = "X"

$mslists.Range("A9").Value = "filter"
$mslists.Range("C9").Value = "X"

$mslists.Range("A10").Value = "groupNames"
$mslists.Range("B10").Value = "X"

$mslists.Range("A11").Value = "initialize"
$mslists.Range("C11").Value = "X"

$mslists.Range("A12").Value = "length"
$mslists.Range("B12").Value = "X"

$mslists.Range("A13").Value = "peakLists"
$mslists.Range("C13").Value = "X"

$mslists.Range("A14").Value = "plotSpec"
$mslists.Range("B14").Value = "X?"

$mslists.Range("A15").Value = "show"
$mslists.Range("C15").Value = "X"

# --- styling: column A method names use the monospaced "Fira Code" font --
$firstName = $mslists.Range("A2")
$firstName.Font.Name = "Fira Code"
$firstName.Font.Size = 10
$firstName.Font.Color = 0
$firstName.VerticalAlignment = -4108

$firstName.Copy() | Out-Null
$nameRange = $mslists.Range("A3:A14")
$nameRange.PasteSpecial(-4122) | Out-Null

# "show" row (A15) gets the same font plus a white fill
$firstName.Copy() | Out-Null
$doneRow = $mslists.Range("A15")
$doneRow.PasteSpecial(-4122) | Out-Null
$doneRow.Interior.Color = 16777215
$excel.CutCopyMode = $false

# --- column widths (approximate best-fit) ---------------------------------
$mslists.Columns.Item(1).ColumnWidth = 20.71
$mslists.Columns.Item(3).ColumnWidth = 11.43
$mslists.Columns.Item(4).ColumnWidth = 10.86
$mslists.Columns.Item(5).ColumnWidth = 13.71

# --- page setup ------------------------------------------------------------
$mslists.PageSetup.PaperSize = 9
$mslists.PageSetup.Orientation = 1

# --- selection / active sheet: mslists becomes the active tab ------------
$mslists.Range("D15").Select() | Out-Null
$mslists.Activate() | Out-Null
